# Updates for gas cooler model
# Clarify the "number of rows" labels for the evaporator coil geometry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Value = "number of rows high (tubes per bank)"
$ws.Range("B17").Value = "number of rows deep (numer of banks)"

# Move the active selection from D15 to D14, matching the updated sheet view.
$ws.Range("D14").Select()
